$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update row 11: "9 Months to 31 March 2016" -> "Year to 30 June 2016" with new figures
$ws.Range("A11").Value = "Year to"
$ws.Range("B11").Value = [datetime]"2016-06-30"
$ws.Range("C11").Value = 9.2
$ws.Range("D11").Value = 0.9
$ws.Range("E11").Value = 8.3

# --- Update row 12: Annualised Performance figures
$ws.Range("C12").Value = 15.6
$ws.Range("D12").Value = 4.8
$ws.Range("E12").Value = 10.9

# --- Remove the blank spacer row above "Cumulative Performance" (old row 13),
#     so Cumulative Performance shifts from row 14 up to row 13
$ws.Rows.Item(13).Delete()

# --- Update the (now) row 13: Cumulative Performance figures
$ws.Range("C13").Value = 327
$ws.Range("D13").Value = 59
$ws.Range("E13").Value = 268

# --- Update view state: scroll down one row, move active cell
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("F11").Select()
